# Apply the edits described by the diff: update the date line and the
# multiplication problems in the table.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-08-26 Saturday" "2023-08-27 Sunday"

Replace-Text "61×42=" "85×35="
Replace-Text "28×56=" "94×48="
Replace-Text "65×94=" "49×59="
Replace-Text "24×11=" "97×37="
Replace-Text "47×31=" "47×84="

Replace-Text "91×61=" "18×95="
Replace-Text "41×61=" "61×19="
Replace-Text "72×56=" "15×91="
Replace-Text "66×61=" "67×43="
Replace-Text "65×11=" "84×96="

Replace-Text "94×64=" "94×85="
Replace-Text "29×60=" "13×55="
Replace-Text "28×47=" "54×89="
Replace-Text "31×28=" "58×91="
Replace-Text "46×17=" "29×62="

Replace-Text "65×19=" "44×64="
Replace-Text "59×27=" "79×95="
Replace-Text "33×72=" "99×79="
Replace-Text "78×52=" "28×90="
Replace-Text "68×96=" "22×36="

Replace-Text "53×97=" "79×12="
Replace-Text "68×62=" "93×64="
Replace-Text "78×66=" "60×47="
Replace-Text "42×84=" "19×39="
Replace-Text "68×21=" "66×45="
